function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range('D2') '42.522.55'
Set-TextValue $ws.Range('E2') '  -1.07%  '
Set-TextValue $ws.Range('D3') '2.228.57'
Set-TextValue $ws.Range('E3') '  -0.50%  '
Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  -0.13%  '
Set-TextValue $ws.Range('D5') '112.59'
Set-TextValue $ws.Range('E5') '  -1.08%  '
Set-TextValue $ws.Range('D6') '295.13'
Set-TextValue $ws.Range('E6') '  +9.66%  '
Set-TextValue $ws.Range('D7') '0.627'
Set-TextValue $ws.Range('E7') '  -0.06%  '
Set-TextValue $ws.Range('E8') '  -0.34%  '
Set-TextValue $ws.Range('D9') '0.600'
Set-TextValue $ws.Range('E9') '  -0.99%  '
Set-TextValue $ws.Range('D10') '43.51'
Set-TextValue $ws.Range('E10') '  -5.83%  '
Set-TextValue $ws.Range('E11') '  -1.04%  '
Set-TextValue $ws.Range('D12') '54.28'
Set-TextValue $ws.Range('E12') '  +0.80%  '
Set-TextValue $ws.Range('D13') '8.70'
Set-TextValue $ws.Range('E13') '  -4.84%  '
Set-TextValue $ws.Range('D14') '1.06'
Set-TextValue $ws.Range('E14') '  +21.10%  '
Set-TextValue $ws.Range('E15') '  -1.36%  '
Set-TextValue $ws.Range('D16') '14.99'
Set-TextValue $ws.Range('E16') '  -2.23%  '
Set-TextValue $ws.Range('D17') '2.564.70'
Set-TextValue $ws.Range('E17') '  -0.52%  '
Set-TextValue $ws.Range('D18') '2.230.08'
Set-TextValue $ws.Range('E18') '  -0.49%  '
Set-TextValue $ws.Range('D19') '42.502.58'
Set-TextValue $ws.Range('E19') '  -0.99%  '
Set-TextValue $ws.Range('D20') '7.23'
Set-TextValue $ws.Range('E20') '  +7.20%  '
Set-TextValue $ws.Range('E21') '  -1.62%  '
Set-TextValue $ws.Range('D22') '73.59'
Set-TextValue $ws.Range('E22') '  +2.19%  '
Set-TextValue $ws.Range('E23') '  +15.76%  '
Set-TextValue $ws.Range('E24') '  +0.48%  '
Set-TextValue $ws.Range('D25') '239.71'
Set-TextValue $ws.Range('E25') '  +3.97%  '
Set-TextValue $ws.Range('D26') '8.90'
Set-TextValue $ws.Range('E26') '  -4.91%  '
Set-TextValue $ws.Range('D27') '1.00'
Set-TextValue $ws.Range('E27') '  -1.20%  '
Set-TextValue $ws.Range('D28') '11.45'
Set-TextValue $ws.Range('E28') '  -6.31%  '
Set-TextValue $ws.Range('B29') 'Monero'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D29') '175.66'
Set-TextValue $ws.Range('E29') '  +1.12%  '
Set-TextValue $ws.Range('B30') 'InjectiveProtocol'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D30') '37.08'
Set-TextValue $ws.Range('E30') '  -8.32%  '
Set-TextValue $ws.Range('B31') 'Toncoin'
Set-TextValue $ws.Range('C31') 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Range('D31') '2.13'
Set-TextValue $ws.Range('E31') '  -4.20%  '
Set-TextValue $ws.Range('D32') '21.63'
Set-TextValue $ws.Range('E32') '  +2.34%  '
Set-TextValue $ws.Range('E33') '  -4.64%  '
Set-TextValue $ws.Range('D34') '0.0881'
Set-TextValue $ws.Range('E34') '  -2.48%  '
Set-TextValue $ws.Range('D35') '5.70'
Set-TextValue $ws.Range('E35') '  +2.23%  '
Set-TextValue $ws.Range('D36') '4.92'
Set-TextValue $ws.Range('E36') '  +4.94%  '
Set-TextValue $ws.Range('E37') '  -0.74%  '
Set-TextValue $ws.Range('E38') '  -2.48%  '
Set-TextValue $ws.Range('D39') '0.0374'
Set-TextValue $ws.Range('E39') '  -0.02%  '
Set-TextValue $ws.Range('E40') '  -2.11%  '
Set-TextValue $ws.Range('D41') '2.40'
Set-TextValue $ws.Range('E41') '  -5.99%  '
Set-TextValue $ws.Range('D42') '71.38'
Set-TextValue $ws.Range('E42') '  +0.48%  '
Set-TextValue $ws.Range('E43') '  -1.80%  '
Set-TextValue $ws.Range('E44') '  +0.19%  '
Set-TextValue $ws.Range('D45') '12.30'
Set-TextValue $ws.Range('E45') '  -6.81%  '
Set-TextValue $ws.Range('E46') '  -2.16%  '
Set-TextValue $ws.Range('D47') '5.42'
Set-TextValue $ws.Range('E47') '  -4.46%  '
Set-TextValue $ws.Range('D48') '1.29'
Set-TextValue $ws.Range('E48') '  +2.90%  '
Set-TextValue $ws.Range('D49') '8.53'
Set-TextValue $ws.Range('E49') '  +0.83%  '
Set-TextValue $ws.Range('D50') '102.24'
Set-TextValue $ws.Range('E50') '  +1.95%  '
Set-TextValue $ws.Range('B51') 'Cronos'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range('D51') '0.0976'
Set-TextValue $ws.Range('E51') '  -1.29%  '
